# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

# "Forecast Comparison" sheet: MyForecast value for week W5 (row 6) updated
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D6").Value = 7

# "Summary" sheet: recomputed totals / forecast week dates.
# These cells hold plain text (e.g. "109", "2025-02-23"), not real numbers
# or dates, so force the Text number format before writing the new value -
# otherwise Excel's automatic type detection would turn "108" into a number
# and "2025-03-09" into a date serial. Resetting the style back to "Normal"
# afterwards keeps the cell's style index identical to the original
# (avoids leaving a stray "@" text format applied to the cell).
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "108"
$wsSummary.Range("B9").Style = "Normal"

$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B13").Value = "2025-03-09"
$wsSummary.Range("B13").Style = "Normal"

$wsSummary.Range("B15").NumberFormat = "@"
$wsSummary.Range("B15").Value = "2025-05-11"
$wsSummary.Range("B15").Style = "Normal"
